# Traveling Salesman - "Added Functionality to TS GA"
# Build an O/Q/S "order check" block (columns O, Q, S) next to the existing
# x/y point table, plus a U column that verifies two computed tours (O vs Q,
# and Q vs S) line up city-for-city.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O: first tour order labels -------------------------------------------------
$oValues = @(
    "C10(100,160)",
    "C17(100,40)",
    "C18(200,160)",
    "C11(120,80)",
    "C19(200,40)",
    "C16(180,100)",
    "C13(180,60)",
    "C9(20,160)",
    "C8(140,140)",
    "C14(180,200)",
    "C15(160,20)",
    "C12(140,180)",
    "C3(40,120)",
    "C2(20,40)",
    "C1(60,20)",
    "C0(20,20)",
    "C7(80,180)",
    "C6(100,120)",
    "C5(60,200)",
    "C4(60,80)"
)

for ($i = 0; $i -lt $oValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 15).Value = $oValues[$i]
}

# --- Column Q: second tour order labels (identical to column O) ------------------------
for ($i = 0; $i -lt $oValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 17).Value = $oValues[$i]
}

# --- Column S: third tour order labels (the one being checked against O/Q) -------------
$sValues = @(
    "C10(100,160)",
    "C17(200,160)",
    "C18(120,80)",
    "C11(200,40)",
    "C19(180,100)",
    "C16(180,60)",
    "C13(100,40)",
    "C9(140,140)",
    "C8(180,200)",
    "C14(20,160)",
    "C15(160,20)",
    "C12(140,180)",
    "C3(20,40)",
    "C2(40,120)",
    "C1(20,20)",
    "C0(80,180)",
    "C7(100,120)",
    "C6(60,20)",
    "C5(60,80)",
    "C4(60,200)"
)

for ($i = 0; $i -lt $sValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 19).Value = $sValues[$i]
}

# --- Column U: boolean "do all three line up on this row?" -----------------------------
# Row 1 gets its own (non-shared) formula.
$ws.Range("U1").Formula = "=AND(O1 = Q1, Q1 = S1)"

# Rows 2-20 are filled as one shared formula (mirrors Excel's fill-down behaviour).
$ws.Range("U2:U20").Formula = "=AND(O2 = Q2, Q2 = S2)"

# --- Put the selection on U1, matching the saved view ----------------------------------
$ws.Range("U1").Select() | Out-Null
